$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Turn off track-changes recording so our own edits are applied as
#    plain edits (not wrapped in new ins/del revisions), then accept
#    every existing tracked change (this is what "Rettet ... efter
#    review" does - the author reviewed and accepted the previous
#    editor's suggestions).
# ------------------------------------------------------------------
$d.TrackRevisions = $false
$d.Revisions.AcceptAll()

# ------------------------------------------------------------------
# 2. Resolve the reviewer comment by rewriting the commented sentence
#    and removing the comment itself.
# ------------------------------------------------------------------
$d.Comments.Item(1).Delete()

$rng1 = $d.Content
$find1 = ", hvis en klasse er udviklet af nogle og den anden klasse som den skal interagere"
$replace1 = ", hvis udviklerne har forstået interaktionen mellem klasserne på forskellige måde"
$rng1.Find.Execute($find1, $true, $false, $false, $false, $false, $true, 1, $false, $replace1, 2) | Out-Null

$rng2 = $d.Content
$find2 = " med har forstået at det fungere på en anden måde"
$rng2.Find.Execute($find2, $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# Move the _GoBack bookmark to mark this as the last edited spot
# (right before the final "r" of "interagerer"/"måder").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$rng3 = $d.Content
$rng3.Find.Execute("måder", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$goBackPos = $rng3.End - 1
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# ------------------------------------------------------------------
# 3. Convert the "Figur" caption's SEQ field (begin/instrText/
#    separate/end fldChar quartet) into a simple field (fldSimple),
#    keeping the bookmark _Ref451344130 wrapped around "Figur <n>".
# ------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_Ref451344130")
$bmStart = $bm.Range.Start

$seqField = $null
for ($i = 1; $i -le $d.Fields.Count; $i++) {
    $candidate = $d.Fields.Item($i)
    if ($candidate.Code.Text -match "SEQ Figur") {
        $seqField = $candidate
    }
}
$seqField.Delete()

$figRng = $d.Content
$figRng.Find.Execute("Figur :", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fieldInsertPos = $figRng.End - 1
$fieldInsertRange = $d.Range($fieldInsertPos, $fieldInsertPos)
$d.Fields.Add($fieldInsertRange, 12, "SEQ Figur \* ARABIC", $false) | Out-Null

$newSeqField = $d.Fields.Item($d.Fields.Count)
$newFieldEnd = $newSeqField.Result.End

$d.Bookmarks.Item("_Ref451344130").Delete()
$newBmRange = $d.Range($bmStart, $newFieldEnd)
$d.Bookmarks.Add("_Ref451344130", $newBmRange) | Out-Null
